$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow a set of columns (raw OOXML width 8 -> 7) ---
# Empirically, setting .ColumnWidth = N results in a stored OOXML width of N + 0.8333333333333333
# so to land on an integer raw width of 7 we use 6.166666666666667
$narrowCols = @(3,5,6,7,9,12,15,16,17,21,22,24,26,29,30,34)
foreach ($c in $narrowCols) {
    $ws.Columns.Item($c).ColumnWidth = 6.166666666666667
}

# Column T (20) goes from raw width 9 -> 8
$ws.Columns.Item(20).ColumnWidth = 7.166666666666667

# --- Update row 5 values to 2 decimal-place "custom accuracy" rounding ---
$ws.Range("C5").Value = 13.81
$ws.Range("D5").Value = 0.83
$ws.Range("E5").Value = 40.05
$ws.Range("F5").Value = 32.64
$ws.Range("G5").Value = 14.04
$ws.Range("H5").Value = 51.02
$ws.Range("I5").Value = 21.87
$ws.Range("J5").Value = 9.68
$ws.Range("K5").Value = 14.54
$ws.Range("L5").Value = 15.93
$ws.Range("M5").Value = 17.39
$ws.Range("N5").Value = 4.59
$ws.Range("O5").Value = 13.87
$ws.Range("P5").Value = 20.39
$ws.Range("Q5").Value = 11.8
$ws.Range("R5").Value = 0.44
$ws.Range("S5").Value = 0.45
$ws.Range("T5").Value = 209.93
$ws.Range("U5").Value = 39.82
$ws.Range("V5").Value = 13.38
$ws.Range("W5").Value = 27.12
$ws.Range("X5").Value = 14.36
$ws.Range("Y5").Value = 1.87
$ws.Range("Z5").Value = 25.38
$ws.Range("AA5").Value = 11.64
$ws.Range("AB5").Value = 9.94
$ws.Range("AC5").Value = 12.12
$ws.Range("AD5").Value = 16.87
$ws.Range("AE5").Value = 0.12
$ws.Range("AF5").Value = 45.61
$ws.Range("AG5").Value = 7.36
$ws.Range("AH5").Value = 16.31

# --- Remove row 6 entirely (data reduced by one record) ---
$ws.Rows.Item(6).Delete()
